# Updated cryptos list - apply new Price (D) and Volume(1h) (E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = "D"; Value = "27.319.62" },
    @{ Row = 2; Col = "E"; Value = "  +0.97%  " },
    @{ Row = 3; Col = "D"; Value = "1.778.57" },
    @{ Row = 3; Col = "E"; Value = "  +4.06%  " },
    @{ Row = 4; Col = "D"; Value = "0.9996" },
    @{ Row = 4; Col = "E"; Value = "  -0.18%  " },
    @{ Row = 5; Col = "D"; Value = "313.96" },
    @{ Row = 5; Col = "E"; Value = "  +1.94%  " },
    @{ Row = 6; Col = "E"; Value = "  -0.12%  " },
    @{ Row = 7; Col = "D"; Value = "0.5244" },
    @{ Row = 7; Col = "E"; Value = "  +10.90%  " },
    @{ Row = 8; Col = "D"; Value = "0.3606" },
    @{ Row = 8; Col = "E"; Value = "  +5.52%  " },
    @{ Row = 9; Col = "D"; Value = "42.57" },
    @{ Row = 9; Col = "E"; Value = "  +1.11%  " },
    @{ Row = 10; Col = "D"; Value = "0.07382" },
    @{ Row = 10; Col = "E"; Value = "  +1.71%  " },
    @{ Row = 11; Col = "D"; Value = "1.093" },
    @{ Row = 11; Col = "E"; Value = "  +5.97%  " },
    @{ Row = 12; Col = "D"; Value = "0.9997" },
    @{ Row = 12; Col = "E"; Value = "  -0.10%  " },
    @{ Row = 13; Col = "E"; Value = "  +4.36%  " },
    @{ Row = 14; Col = "E"; Value = "  +4.21%  " },
    @{ Row = 15; Col = "D"; Value = "1.779.47" },
    @{ Row = 15; Col = "E"; Value = "  +4.24%  " },
    @{ Row = 16; Col = "D"; Value = "6.995" },
    @{ Row = 16; Col = "E"; Value = "  +2.58%  " },
    @{ Row = 17; Col = "D"; Value = "88.49" },
    @{ Row = 18; Col = "D"; Value = "0.00001046" },
    @{ Row = 18; Col = "E"; Value = "  +1.18%  " },
    @{ Row = 19; Col = "E"; Value = "  +1.03%  " },
    @{ Row = 20; Col = "D"; Value = "0.9993" },
    @{ Row = 20; Col = "E"; Value = "  -0.13%  " },
    @{ Row = 21; Col = "D"; Value = "16.75" },
    @{ Row = 21; Col = "E"; Value = "  +1.92%  " },
    @{ Row = 22; Col = "D"; Value = "5.850" },
    @{ Row = 22; Col = "E"; Value = "  +4.67%  " },
    @{ Row = 23; Col = "D"; Value = "27.402.06" },
    @{ Row = 23; Col = "E"; Value = "  +1.17%  " },
    @{ Row = 24; Col = "D"; Value = "11.33" },
    @{ Row = 24; Col = "E"; Value = "  +4.69%  " },
    @{ Row = 25; Col = "D"; Value = "2.069" },
    @{ Row = 25; Col = "E"; Value = "  -2.02%  " },
    @{ Row = 26; Col = "D"; Value = "153.51" },
    @{ Row = 26; Col = "E"; Value = "  -1.91%  " },
    @{ Row = 27; Col = "D"; Value = "20.10" },
    @{ Row = 27; Col = "E"; Value = "  +3.00%  " },
    @{ Row = 28; Col = "D"; Value = "2.349" },
    @{ Row = 28; Col = "E"; Value = "  +14.01%  " },
    @{ Row = 29; Col = "D"; Value = "1.986.19" },
    @{ Row = 29; Col = "E"; Value = "  +4.43%  " },
    @{ Row = 30; Col = "D"; Value = "121.47" },
    @{ Row = 30; Col = "E"; Value = "  +2.01%  " },
    @{ Row = 31; Col = "D"; Value = "1.063" },
    @{ Row = 31; Col = "E"; Value = "  +5.67%  " },
    @{ Row = 32; Col = "D"; Value = "0.09787" },
    @{ Row = 32; Col = "E"; Value = "  +7.09%  " },
    @{ Row = 33; Col = "D"; Value = "5.545" },
    @{ Row = 33; Col = "E"; Value = "  +4.99%  " },
    @{ Row = 34; Col = "D"; Value = "3.604" },
    @{ Row = 34; Col = "E"; Value = "  +0.68%  " },
    @{ Row = 35; Col = "D"; Value = "0.02238" },
    @{ Row = 35; Col = "E"; Value = "  +2.64%  " },
    @{ Row = 36; Col = "D"; Value = "0.05990" },
    @{ Row = 36; Col = "E"; Value = "  +3.30%  " },
    @{ Row = 37; Col = "D"; Value = "11.23" },
    @{ Row = 37; Col = "E"; Value = "  +2.26%  " },
    @{ Row = 38; Col = "D"; Value = "4.855" },
    @{ Row = 38; Col = "E"; Value = "  +2.91%  " },
    @{ Row = 39; Col = "D"; Value = "0.6158" },
    @{ Row = 39; Col = "E"; Value = "  +5.09%  " },
    @{ Row = 40; Col = "D"; Value = "0.2029" },
    @{ Row = 40; Col = "E"; Value = "  +2.51%  " },
    @{ Row = 41; Col = "D"; Value = "1.426" },
    @{ Row = 41; Col = "E"; Value = "  +2.61%  " },
    @{ Row = 42; Col = "D"; Value = "8.087" },
    @{ Row = 42; Col = "E"; Value = "  +8.68%  " },
    @{ Row = 43; Col = "D"; Value = "1.149" },
    @{ Row = 43; Col = "E"; Value = "  +4.82%  " },
    @{ Row = 44; Col = "E"; Value = "  +5.06%  " },
    @{ Row = 45; Col = "D"; Value = "0.5772" },
    @{ Row = 45; Col = "E"; Value = "  +2.81%  " },
    @{ Row = 46; Col = "D"; Value = "3.634" },
    @{ Row = 46; Col = "E"; Value = "  +2.17%  " },
    @{ Row = 47; Col = "D"; Value = "121.53" },
    @{ Row = 47; Col = "E"; Value = "  +3.44%  " },
    @{ Row = 48; Col = "D"; Value = "1.892" },
    @{ Row = 48; Col = "E"; Value = "  +3.50%  " },
    @{ Row = 49; Col = "D"; Value = "1.111" },
    @{ Row = 49; Col = "E"; Value = "  +2.91%  " },
    @{ Row = 50; Col = "D"; Value = "0.06722" },
    @{ Row = 50; Col = "E"; Value = "  +1.49%  " },
    @{ Row = 51; Col = "D"; Value = "70.96" },
    @{ Row = 51; Col = "E"; Value = "  +2.13%  " }
)

foreach ($u in $updates) {
    $addr = "$($u.Col)$($u.Row)"
    $rng = $ws.Range($addr)
    if ($u.Col -eq "D") {
        # Price column sometimes looks like a plain number (e.g. "313.96");
        # force text so Excel doesn't silently reinterpret it as a Number.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
